# Log of Meetings.xlsx — add two new meeting-minute rows (5 & 6) to Sheet1
# and move the active-cell selection, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 : 9/27/2023 meeting -------------------------------------------
# Column A: date, formatted like the rows above it (copy style from A2).
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 45196

# Column C (Time Ended) entered before column B so the shared-string table
# picks up "4:15pm" ahead of "3:45pm", matching the authored file.
$ws.Range("C5").Value = "4:15pm"

# Column B (Time Started): was originally keyed in as an actual time value
# (which stamps the cell with an h:mm number format) and later overwritten
# with literal text - reproduce that so the cell keeps the h:mm style.
$ws.Range("B5").Value = 0.65625
$ws.Range("B5").NumberFormat = "h:mm"
$ws.Range("B5").Value = "3:45pm"

$ws.Range("D5").Value = "30 minutes"
$ws.Range("E5").Value = "Progress update and discussion of libraries for skeletal extraction. Set up a meeting with client"

# --- Row 6 : 10/6/2023 meeting --------------------------------------------
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 45205

$ws.Range("B6").Value = "3:00pm"
$ws.Range("C6").Value = "3:30pm"
$ws.Range("D6").Value = "30 minutes"
$ws.Range("E6").Value = "Client meeting to discuss progress. Was given access to a layout for our project and videos for analysis. Switched first exercise to deep squat"

# --- Selection -------------------------------------------------------------
$ws.Range("D11").Select() | Out-Null
